$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.953.14"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "1.654.48"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.511"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0875"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").Value = "1.888.17"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "1.645.58"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "26.956.61"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").Value = "0.0₃0732"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "1.542.97"
$ws.Range("E32").Value = "  +3.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("E34").Value = "  +4.81%  "
$ws.Range("E35").Value = "  +8.45%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.583"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.891"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.43%  "
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.88%  "
$ws.Range("D44").Value = "1.795.40"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.774"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.917"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0989"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0506"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
